$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (shifts rows 48..147 down to 49..148)
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the data from the diff
$ws.Range("A48").Value = 11
$ws.Range("B48").Value = "Vega Monumental Concepción"
$ws.Range("C48").Value = "Bíobío"
$ws.Range("D48").Value = 45238
$ws.Range("E48").Value = 8
$ws.Range("F48").Value = 100112037
$ws.Range("G48").Value = "Cebollín"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 100
$ws.Range("K48").Value = 3200
$ws.Range("L48").Value = 3400
$ws.Range("M48").Value = 3300
$ws.Range("N48").Value = "$/paquete 36 unidades"
$ws.Range("O48").Value = "Región Metropolitana"
$ws.Range("P48").Value = 92
$ws.Range("Q48").Value = 36
$ws.Range("R48").Value = "Hortaliza"
